# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 04:49:32"
$wsZhCn.Range("H2").Value = "2016-03-12 04:49:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 04:49:37"
$wsDeDe.Range("H2").Value = "2016-03-12 04:49:53"
